$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 22.24000000000004
$ws.Range("G2").Value = 0.00000000000002897682094271659
$ws.Range("H2").Value = 0.00000000000009480952178544988
$ws.Range("K2").Value = 45.01736723312373
$ws.Range("L2").Value = "[33.651238389443044, 56.38349607680441]"
$ws.Range("M2").Value = 0.0000000000003770317391627032
$ws.Range("N2").Value = 0.0000000000003770317391627032
$ws.Range("O2").Value = 1.553500271144502
$ws.Range("P2").Value = "[1.2641844311742716, 1.8428161111147325]"
$ws.Range("S2").Value = 57.67545599889662
$ws.Range("T2").Value = "[50.33796314019112, 65.01294885760211]"
$ws.Range("W2").Value = 16.74122122122125
$ws.Range("X2").Value = 15.71715715715719
$ws.Range("Y2").Value = 17.76528528528532

# Row 3 updates
$ws.Range("E3").Value = 24.78000000000043
$ws.Range("H3").Value = 0.0000000000000004249657510526915
$ws.Range("K3").Value = 47.60801289875827
$ws.Range("L3").Value = "[36.89626827444492, 58.319757523071615]"
$ws.Range("O3").Value = -2.188737224122619
$ws.Range("P3").Value = "[-2.415158316273234, -1.9623161319720035]"
$ws.Range("S3").Value = 62.80621707082722
$ws.Range("T3").Value = "[57.261934410892025, 68.35049973076241]"
$ws.Range("W3").Value = 8.632072072072223
$ws.Range("X3").Value = 7.739099099099236
$ws.Range("Y3").Value = 9.525045045045211
